# Apply the commit: insert a new client "GONZALEZ CARDENAS ERNESTO PAOLO"
# (with zero sales) into the "LOZANO MOLINA TITO" salesperson group, which
# is alphabetically sorted by client name (column B) within each
# salesperson group (column A) on both the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets. The new row pushes every following row down by
# one, and the "x de N" summary counts on the VENTAS POR GRUPO sheet are
# refreshed to reflect the new total row count (323 -> 324).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": 16 data columns (C..R), new row at 244
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(244).Insert()
$ws1.Range("A244").Value = "LOZANO MOLINA TITO"
$ws1.Range("B244").Value = "GONZALEZ CARDENAS ERNESTO PAOLO"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(244, $col).Value = 0
}

# Refresh the "x de 323" -> "x de 324" summary row (now row 326)
$ws1.Range("C326").Value = "5 de 324"
$ws1.Range("D326").Value = "14 de 324"
$ws1.Range("E326").Value = "10 de 324"
$ws1.Range("F326").Value = "0 de 324"
$ws1.Range("G326").Value = "0 de 324"
$ws1.Range("H326").Value = "10 de 324"
$ws1.Range("I326").Value = "15 de 324"
$ws1.Range("J326").Value = "0 de 324"
$ws1.Range("K326").Value = "4 de 324"
$ws1.Range("L326").Value = "22 de 324"
$ws1.Range("M326").Value = "49 de 324"
$ws1.Range("N326").Value = "2 de 324"
$ws1.Range("O326").Value = "1 de 324"
$ws1.Range("P326").Value = "1 de 324"
$ws1.Range("Q326").Value = "0 de 324"
$ws1.Range("R326").Value = "0 de 324"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": 5 data columns (C..G), new row at 248
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(248).Insert()
$ws2.Range("A248").Value = "LOZANO MOLINA TITO"
$ws2.Range("B248").Value = "GONZALEZ CARDENAS ERNESTO PAOLO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(248, $col).Value = 0
}
